$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 data entry (melhorias na leitura de dados)
$ws.Range("A3").Value = "10/07/2023 15:31"

$ws.Range("B3").Value = 0.25
$ws.Range("C3").Formula = "=B3 * 600 / 100 *100"

$ws.Range("D3").Value = 0.8
$ws.Range("E3").Formula = "=D3 * 600 / 100 *100"

$ws.Range("F3").Value = 0.4
$ws.Range("G3").Formula = "=F3 * 600 / 100 *100"

$ws.Range("H3").Value = 0.16
$ws.Range("I3").Formula = "=H3 * 600 / 100 *100"

$ws.Range("J3").Value = 0.65
$ws.Range("K3").Formula = "=J3 * 600 / 100 *100"

$ws.Range("L3").Value = 0.9399999999999999
$ws.Range("M3").Formula = "=L3 * 600 / 100 *100"

$ws.Range("N3").Formula = "=C3+E3+G3+I3+K3+M3"
